$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text (8/28/19 -> 10/26/19)
#    on the slide master and on every slide layout (12 occurrences total).
# ---------------------------------------------------------------------------
$newDate = "10/26/19"

function Set-DatePlaceholderText($shapes, $text) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.Designs.Item(1).SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# ---------------------------------------------------------------------------
# 2) Slide 1 subtitle: drop the "India ML Hiring Hackathon 2019" paragraph,
#    keep only "Sachin Rastogi".
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "Sachin Rastogi"

# ---------------------------------------------------------------------------
# 3) Slide 8 title text change.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$title8 = $slide8.Shapes.Item(1)
$title8.TextFrame.TextRange.Text = "5 things to focus on while solving such problems?"
